# Auto-generated Excel COM-interop script applying the diff from commit
# 'chore: update Sheets via scheduled runner' to Ravana_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1626.4445
$ws.Range("I9").Value = 1850.25
$ws.Range("J9").Value = 1447.4
$ws.Range("K9").Value = 1850.25
$ws.Range("L9").Value = 1447.4
$ws.Range("M9").Value = -1681.25
$ws.Range("N9").Value = -1785.4
$ws.Range("H11").Value = 169.25
$ws.Range("I11").Value = 169.25
$ws.Range("K11").Value = 169.25
$ws.Range("M11").Value = -29.25
$ws.Range("H40").Value = 835
$ws.Range("I40").Value = 820
$ws.Range("K40").Value = 820
$ws.Range("M40").Value = -645
$ws.Range("H98").Value = 1648.5
$ws.Range("J98").Value = 1999.5
$ws.Range("L98").Value = 1999.5
$ws.Range("N98").Value = -4995.5
$ws.Range("H122").Value = 1648.5
$ws.Range("J122").Value = 1999.5
$ws.Range("L122").Value = 5998.5
$ws.Range("N122").Value = -10898.5
$ws.Range("H137").Value = 4677.1
$ws.Range("I137").Value = 2663.4285
$ws.Range("K137").Value = 7990.2855
$ws.Range("M137").Value = -5440.2855
$ws.Range("H138").Value = 4335.613
$ws.Range("I138").Value = 3066
$ws.Range("J138").Value = 4579.769
$ws.Range("K138").Value = 9198
$ws.Range("L138").Value = 13739.307
$ws.Range("M138").Value = -4058
$ws.Range("N138").Value = -24019.307
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 780.5
$ws.Range("I2").Value = 744.1429000000001
$ws.Range("K2").Value = 744.1429000000001
$ws.Range("M2").Value = -631.1429000000001
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 950
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 950
$ws.Range("M4").Value = 66
$ws.Range("N4").Value = -1182
$ws.Range("H32").Value = 3762.318
$ws.Range("I32").Value = 2388.6
$ws.Range("K32").Value = 2388.6
$ws.Range("M32").Value = -2101.6
$ws.Range("H61").Value = 3774.5
$ws.Range("I61").Value = 3774.5
$ws.Range("K61").Value = 3774.5
$ws.Range("M61").Value = -3562.5
$ws.Range("H116").Value = 780.5
$ws.Range("I116").Value = 744.1429000000001
$ws.Range("K116").Value = 744.1429000000001
$ws.Range("M116").Value = 1549.8571
$ws.Range("H132").Value = 4161.091
$ws.Range("I132").Value = 2967.7144
$ws.Range("J132").Value = 6249.5
$ws.Range("K132").Value = 8903.143199999999
$ws.Range("L132").Value = 18748.5
$ws.Range("M132").Value = -6373.143199999999
$ws.Range("N132").Value = -23808.5
$ws.Range("H136").Value = 3774.5
$ws.Range("I136").Value = 3774.5
$ws.Range("K136").Value = 11323.5
$ws.Range("M136").Value = -8773.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 780.5
$ws.Range("I3").Value = 744.1429000000001
$ws.Range("K3").Value = 744.1429000000001
$ws.Range("M3").Value = -630.1429000000001
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H134").Value = 2570.3
$ws.Range("I134").Value = 2570.3
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7710.900000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5175.900000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1866.8823
$ws.Range("J31").Value = 1808.75
$ws.Range("L31").Value = 1808.75
$ws.Range("N31").Value = -2398.75
$ws.Range("H34").Value = 1866.8823
$ws.Range("J34").Value = 1808.75
$ws.Range("L34").Value = 1808.75
$ws.Range("N34").Value = -2212.75
$ws.Range("H132").Value = 5113.125
$ws.Range("I132").Value = 4228
$ws.Range("K132").Value = 12684
$ws.Range("M132").Value = -10154
$ws.Range("H134").Value = 2976.889
$ws.Range("I134").Value = 2976.889
$ws.Range("K134").Value = 8930.667000000001
$ws.Range("M134").Value = -6395.667000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 467.36365
$ws.Range("I12").Value = 419.4
$ws.Range("J12").Value = 507.33334
$ws.Range("K12").Value = 1258.2
$ws.Range("L12").Value = 1522.00002
$ws.Range("M12").Value = -1085.2
$ws.Range("N12").Value = -1868.00002
$ws.Range("H33").Value = 187.25
$ws.Range("I33").Value = 99
$ws.Range("J33").Value = 216.66667
$ws.Range("K33").Value = 594
$ws.Range("L33").Value = 1300.00002
$ws.Range("M33").Value = -311
$ws.Range("N33").Value = -1866.00002
$ws.Range("H34").Value = 2627.1428
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 2981.6667
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 8945.000100000001
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -9113.000100000001
$ws.Range("H55").Value = 2750
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2750
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 8250
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -8604
$ws.Range("H80").Value = 6164
$ws.Range("I80").Value = 4751
$ws.Range("J80").Value = 8990
$ws.Range("K80").Value = 14253
$ws.Range("L80").Value = 26970
$ws.Range("M80").Value = -13317
$ws.Range("N80").Value = -28842
$ws.Range("H83").Value = 6164
$ws.Range("I83").Value = 4751
$ws.Range("J83").Value = 8990
$ws.Range("K83").Value = 42759
$ws.Range("L83").Value = 80910
$ws.Range("M83").Value = -38079
$ws.Range("N83").Value = -90270
$ws.Range("H98").Value = 3628
$ws.Range("I98").Value = 880
$ws.Range("K98").Value = 2640
$ws.Range("M98").Value = -1142
$ws.Range("H131").Value = 1190.8572
$ws.Range("J131").Value = 1342.75
$ws.Range("L131").Value = 4028.25
$ws.Range("N131").Value = -14108.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 99.666664
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 100
$ws.Range("N2").Value = -326
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3869.1428
$ws.Range("J7").Value = 5999.5
$ws.Range("L7").Value = 5999.5
$ws.Range("N7").Value = -6223.5
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H126").Value = 3869.1428
$ws.Range("J126").Value = 5999.5
$ws.Range("L126").Value = 17998.5
$ws.Range("N126").Value = -22938.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H107").Value = 433
$ws.Range("I107").Value = 422.16666
$ws.Range("K107").Value = 1266.49998
$ws.Range("M107").Value = 653.5000199999999
$ws.Range("H126").Value = 1161.6
$ws.Range("I126").Value = 985.4167
$ws.Range("K126").Value = 2956.2501
$ws.Range("M126").Value = -486.2501000000002
$ws.Range("H132").Value = 2475.8667
$ws.Range("I132").Value = 1713.4445
$ws.Range("J132").Value = 3619.5
$ws.Range("K132").Value = 5140.333500000001
$ws.Range("L132").Value = 10858.5
$ws.Range("M132").Value = -2610.333500000001
$ws.Range("N132").Value = -15918.5
